$d = $word.ActiveDocument

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
#
# header1.xml  -> Headers.Item(2) (first-page header) -> BTec_Logo-Orange, image1.jpg -> image2.jpg
# footer1.xml  -> Footers.Item(2) (first-page footer) -> PearsonLogo, id=3, image2.png -> image1.png
# footer2.xml  -> Footers.Item(1) (default footer)    -> PearsonLogo, id=2, image2.png -> image1.png

$firstHeader = $d.Sections.First.Headers.Item(2)
if ($firstHeader.Exists) {
    $btecShape = $firstHeader.Range.InlineShapes.Item(1)
    $btecShape.Name = "image2.jpg"
}

$firstFooter = $d.Sections.First.Footers.Item(2)
if ($firstFooter.Exists) {
    $pearsonShape1 = $firstFooter.Range.InlineShapes.Item(1)
    $pearsonShape1.Name = "image1.png"
}

$defaultFooter = $d.Sections.First.Footers.Item(1)
if ($defaultFooter.Exists) {
    $pearsonShape2 = $defaultFooter.Range.InlineShapes.Item(1)
    $pearsonShape2.Name = "image1.png"
}
